# Document the new "allowsplits" parameter on the "parameter" sheet.
# Default changed to false (per commit message); the example sheet still
# illustrates the (non-default) "true" setting as text, matching how the
# other parameter rows store their sample value as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameter")
$ws.Activate() | Out-Null

$ws.Range("A5").Value = "allowsplits"

# Force the sample value into the cell as literal text "true" (not a
# boolean) - matches the existing "days"/"9" style sample values in this
# column. Enter it as a formula returning the text, then bake it down to
# a static value via paste-special so it lands as a plain string cell.
$ws.Range("B5").Formula = "=""true"""
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4163) | Out-Null

$ws.Range("C5").Value = "Controls whether a sales order or forecast can be split across multiple manufacturing orders during planning. Default: false"

# Selecting the whole new row leaves the workbook with this sheet active,
# matching the saved selection/active-tab state of the edited file.
$ws.Rows.Item(5).Select() | Out-Null
